$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$text)
    $escaped = $text.Replace("""", """""")
    $cellRange.Formula = "=""" + $escaped + """"
    $cellRange.Copy()
    $cellRange.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

$ws.Range("D2").Value = '67.419.05'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '3.489.86'
$ws.Range("E3").Value = '  -0.94%  '
Set-TextValue $ws.Range("D5") '598.07'
$ws.Range("E5").Value = '  +0.25%  '
Set-TextValue $ws.Range("D6") '178.34'
$ws.Range("E6").Value = '  +2.85%  '
Set-TextValue $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.44%  '
Set-TextValue $ws.Range("D9") '0.134'
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("E10").Value = '  -2.85%  '
Set-TextValue $ws.Range("D11") '0.429'
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("D12").Value = '4.092.78'
$ws.Range("E12").Value = '  -0.91%  '
Set-TextValue $ws.Range("D13") '32.03'
$ws.Range("E13").Value = '  +11.48%  '
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").Value = '67.346.06'
$ws.Range("E15").Value = '  +0.06%  '
Set-TextValue $ws.Range("D16") '0.0000178'
$ws.Range("E16").Value = '  -2.31%  '
$ws.Range("D17").Value = '3.479.94'
$ws.Range("E17").Value = '  -0.88%  '
Set-TextValue $ws.Range("D18") '6.29'
$ws.Range("E18").Value = '  -0.92%  '
Set-TextValue $ws.Range("D19") '14.55'
$ws.Range("E19").Value = '  +1.77%  '
Set-TextValue $ws.Range("D20") '390.73'
$ws.Range("E20").Value = '  -1.94%  '
Set-TextValue $ws.Range("D21") '7.96'
$ws.Range("E21").Value = '  -0.60%  '
Set-TextValue $ws.Range("D22") '73.13'
$ws.Range("E22").Value = '  -0.47%  '
Set-TextValue $ws.Range("D23") '0.999'
$ws.Range("E23").Value = '  -0.03%  '
Set-TextValue $ws.Range("D24") '0.537'
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("E25").Value = '  +0.40%  '
Set-TextValue $ws.Range("D26") '0.0000122'
$ws.Range("E26").Value = '  -0.81%  '
Set-TextValue $ws.Range("D27") '10.28'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("E28").Value = '  -1.04%  '
Set-TextValue $ws.Range("D29") '0.997'
$ws.Range("E29").Value = '  -0.09%  '
Set-TextValue $ws.Range("D30") '6.20'
$ws.Range("E30").Value = '  -1.53%  '
Set-TextValue $ws.Range("D31") '1.43'
$ws.Range("E31").Value = '  -1.81%  '
$ws.Range("E32").Value = '  -1.64%  '
Set-TextValue $ws.Range("D33") '23.68'
Set-TextValue $ws.Range("D34") '7.32'
$ws.Range("E34").Value = '  -1.10%  '
$ws.Range("E35").Value = '  +0.22%  '
Set-TextValue $ws.Range("D36") '163.59'
$ws.Range("E36").Value = '  -0.27%  '
Set-TextValue $ws.Range("D37") '1.94'
$ws.Range("E37").Value = '  +1.06%  '
Set-TextValue $ws.Range("D38") '0.871'
$ws.Range("E38").Value = '  -2.73%  '
Set-TextValue $ws.Range("D39") '7.03'
$ws.Range("E39").Value = '  +1.51%  '
Set-TextValue $ws.Range("D40") '27.60'
$ws.Range("E40").Value = '  +0.58%  '
Set-TextValue $ws.Range("D41") '4.67'
$ws.Range("E41").Value = '  -1.46%  '
Set-TextValue $ws.Range("D42") '26.52'
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").Value = '2.826.26'
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("E46").Value = '  -1.57%  '
$ws.Range("E49").Value = '  -2.58%  '
Set-TextValue $ws.Range("D50") '33.33'
$ws.Range("E50").Value = '  -1.32%  '
Set-TextValue $ws.Range("D51") '6.42'
$ws.Range("E51").Value = '  -2.01%  '

# Row identity swaps (coin rank changes)
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D44") '2.62'
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D45") '0.0724'
$ws.Range("E45").Value = '  -2.86%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D47") '0.0302'
$ws.Range("E47").Value = '  -1.83%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D48") '341.88'
$ws.Range("E48").Value = '  +0.15%  '
